$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format before assigning, so that numeric-looking
# strings (e.g. "591.66") are not auto-converted to native numbers by Excel,
# preserving the original inline-string / text semantics of the workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.295.33'
$ws.Range('E2').Value = '  -4.96%  '

$ws.Range('D3').Value = '3.269.61'
$ws.Range('E3').Value = '  -7.48%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = '591.66'
$ws.Range('E5').Value = '  -4.73%  '

$ws.Range('D6').Value = '151.36'
$ws.Range('E6').Value = '  -12.38%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('D8').Value = '3.263.32'
$ws.Range('E8').Value = '  -7.55%  '

$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  -10.76%  '

$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -13.98%  '

$ws.Range('D11').Value = '6.60'
$ws.Range('E11').Value = '  -9.03%  '

$ws.Range('D12').Value = '0.510'
$ws.Range('E12').Value = '  -13.01%  '

$ws.Range('D13').Value = '38.59'
$ws.Range('E13').Value = '  -16.52%  '

$ws.Range('D14').Value = '0.0000244'
$ws.Range('E14').Value = '  -11.40%  '

$ws.Range('D15').Value = '3.804.67'
$ws.Range('E15').Value = '  -7.28%  '

$ws.Range('D16').Value = '67.424.41'
$ws.Range('E16').Value = '  -4.92%  '

$ws.Range('D17').Value = '3.279.51'
$ws.Range('E17').Value = '  -7.19%  '

$ws.Range('D18').Value = '536.08'
$ws.Range('E18').Value = '  -11.78%  '

$ws.Range('D19').Value = '7.24'
$ws.Range('E19').Value = '  -14.15%  '

$ws.Range('E20').Value = '  -6.18%  '

$ws.Range('D21').Value = '15.09'
$ws.Range('E21').Value = '  -14.62%  '

$ws.Range('D22').Value = '0.761'
$ws.Range('E22').Value = '  -13.64%  '

$ws.Range('D23').Value = '7.86'
$ws.Range('E23').Value = '  -13.53%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '85.75'
$ws.Range('E24').Value = '  -12.29%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '13.59'
$ws.Range('E25').Value = '  -13.39%  '

$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').Value = '3.25'
$ws.Range('E27').Value = '  -12.05%  '

$ws.Range('D28').Value = '8.10'
$ws.Range('E28').Value = '  -10.94%  '

$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '2.16'
$ws.Range('E29').Value = '  -16.00%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '29.40'
$ws.Range('E30').Value = '  -12.80%  '

$ws.Range('D31').Value = '2.66'
$ws.Range('E31').Value = '  -11.35%  '

$ws.Range('E32').Value = '  -11.00%  '

$ws.Range('D33').Value = '6.60'
$ws.Range('E33').Value = '  -18.86%  '

$ws.Range('D34').Value = '537.08'
$ws.Range('E34').Value = '  -12.87%  '

$ws.Range('D35').Value = '5.79'
$ws.Range('E35').Value = '  -15.28%  '

$ws.Range('E36').Value = '  +0.20%  '

$ws.Range('D37').Value = '0.0457'
$ws.Range('E37').Value = '  -8.67%  '

$ws.Range('D38').Value = '53.35'
$ws.Range('E38').Value = '  -6.37%  '

$ws.Range('D39').Value = '0.0860'
$ws.Range('E39').Value = '  -13.62%  '

$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '9.09'
$ws.Range('E40').Value = '  -16.33%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.128'
$ws.Range('E41').Value = '  -10.94%  '

$ws.Range('D42').Value = '2.74'
$ws.Range('E42').Value = '  -19.10%  '

$ws.Range('D43').Value = '2.936.25'
$ws.Range('E43').Value = '  -12.34%  '

$ws.Range('E44').Value = '  -13.03%  '

$ws.Range('D45').Value = '0.0₃0593'
$ws.Range('E45').Value = '  -18.76%  '

$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  -12.06%  '

$ws.Range('D47').Value = '26.77'
$ws.Range('E47').Value = '  -15.93%  '

$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  -0.09%  '

$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').Value = '2.35'
$ws.Range('E49').Value = '  -18.89%  '

$ws.Range('D50').Value = '125.75'
$ws.Range('E50').Value = '  -6.17%  '

$ws.Range('E51').Value = '  -12.16%  '

# Restore default cell style (removes the explicit number-format style we applied
# above) so the saved cells match the original unstyled inline-string cells.
$ws.Range("D2:E51").Style = "Normal"
